$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on D-column cells whose new value looks numeric,
# so Excel keeps them as literal text (matching the original inline-string cells)
# instead of auto-converting to a Number and dropping things like trailing zeros.
$textFormatRows = @(5,6,7,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,26,27,28,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "20.238.71"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.440.95"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").Value = "0.9134"
$ws.Range("E5").Value = "  -8.09%  "
$ws.Range("D6").Value = "275.35"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").Value = "0.3626"
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").Value = "0.3079"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").Value = "38.93"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").Value = "1.025"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "0.06505"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "0.9984"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "5.347"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").Value = "17.49"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "6.051"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "0.00001009"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "1.438.50"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "0.9314"
$ws.Range("E18").Value = "  -6.27%  "
$ws.Range("D19").Value = "0.05618"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").Value = "67.72"
$ws.Range("E20").Value = "  -5.66%  "
$ws.Range("D21").Value = "5.409"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").Value = "14.23"
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("D23").Value = "10.84"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").Value = "2.240"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "20.234.87"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "138.81"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").Value = "2.129"
$ws.Range("E27").Value = "  -6.59%  "
$ws.Range("D28").Value = "16.89"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "1.590.92"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "110.08"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "3.891"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").Value = "0.8064"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "4.832"
$ws.Range("E33").Value = "  -9.61%  "
$ws.Range("D34").Value = "0.07654"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "1.477"
$ws.Range("D36").Value = "0.05840"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "4.668"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("D38").Value = "1.129"
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("D39").Value = "0.01988"
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").Value = "10.18"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("D41").Value = "0.1851"
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("D42").Value = "0.9269"
$ws.Range("D43").Value = "7.185"
$ws.Range("E43").Value = "  -14.16%  "
$ws.Range("D44").Value = "0.5205"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").Value = "3.489"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "11.79"
$ws.Range("E46").Value = "  -4.45%  "
$ws.Range("D47").Value = "116.55"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").Value = "0.5083"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").Value = "1.734"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").Value = "0.06378"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("D51").Value = "0.9884"
$ws.Range("E51").Value = "  -0.07%  "
